$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 22
$ws.Range("G4").Value = 128
$ws.Range("F5").Value = 984
$ws.Range("F9").Value = 128
$ws.Range("F10").Value = 781
$ws.Range("F11").Value = 824
$ws.Range("F12").Value = 242
$ws.Range("F13").Value = 599
$ws.Range("F14").Value = 1180
$ws.Range("F16").Value = 725
$ws.Range("F17").Value = 648
$ws.Range("F18").Value = 250
$ws.Range("F19").Value = 347
$ws.Range("F20").Value = 337
$ws.Range("F21").Value = 711
$ws.Range("F22").Value = 789
$ws.Range("F23").Value = 8308
$ws.Range("F24").Value = 547
$ws.Range("F25").Value = 547
$ws.Range("F30").Value = 221
$ws.Range("F31").Value = 1705
$ws.Range("F32").Value = 17
$ws.Range("F33").Value = 252
$ws.Range("F35").Value = 165
$ws.Range("F37").Value = 228
$ws.Range("F39").Value = 44
$ws.Range("F42").Value = 164
$ws.Range("F44").Value = 141
$ws.Range("F45").Value = 31
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 95
$ws.Range("F15").Value = 51
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 788
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 788
$ws.Range("F3").Value = 22
$ws.Range("G4").Value = 128
$ws.Range("F5").Value = 984
$ws.Range("F9").Value = 128
$ws.Range("F10").Value = 781
$ws.Range("F13").Value = 824
$ws.Range("F14").Value = 242
$ws.Range("F15").Value = 599
$ws.Range("F16").Value = 1180
$ws.Range("F19").Value = 95
$ws.Range("F20").Value = 725
$ws.Range("F21").Value = 648
$ws.Range("F22").Value = 250
$ws.Range("F23").Value = 347
$ws.Range("F24").Value = 337
$ws.Range("F25").Value = 789
$ws.Range("F26").Value = 8308
$ws.Range("F28").Value = 547
$ws.Range("F29").Value = 547
$ws.Range("F32").Value = 221
$ws.Range("F33").Value = 1705
$ws.Range("F34").Value = 17
$ws.Range("F35").Value = 252
$ws.Range("F37").Value = 165
$ws.Range("F41").Value = 51
$ws.Range("F42").Value = 228
$ws.Range("F45").Value = 44
$ws.Range("F48").Value = 164
$ws.Range("F50").Value = 141

Write-Output "Applied all updates"
